# ---------------------------------------------------------------------------
# Edit: 
#  1) Change the table style (Table Design gallery) of the table on slide 16
#     from "Table_0" ({262DFD78-7A93-4D1D-BDA4-3E003C475046}) to the built-in
#     style {94694250-D4B9-415A-9BCA-244933DA2B70}.
#  2) Re-colour the deck's theme (Design > Colors) from the custom "Integral"
#     palette to the stock "Office" palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style ---------------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{94694250-D4B9-415A-9BCA-244933DA2B70}", $true)
    }
}

# --- 2. Theme colours --------------------------------------------------------
function HexToVbRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme palette (dk1..folHlink), replacing the former "Integral" one.
$officeColors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = HexToVbRgb $officeColors[$i - 1]
}
